# Apply the weekly update: insert two new price rows (218-219) for "Apio"
# (Americana (o), Primera/Segunda, Región de Coquimbo) just above the row
# that used to be 218, shifting the subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 218, pushing old rows 218:233 down to 220:235.
$ws.Rows("218:219").Insert()

# Row 218 - new "Primera" quality entry.
$ws.Range("A218").Value = 9
$ws.Range("B218").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C218").Value = "Metropolitana"
$ws.Range("D218").Value = 44610
$ws.Range("E218").Value = 13
$ws.Range("F218").Value = 100112017
$ws.Range("G218").Value = "Apio"
$ws.Range("H218").Value = "Americana (o)"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 61
$ws.Range("K218").Value = 7000
$ws.Range("L218").Value = 8000
$ws.Range("M218").Value = 7508
$ws.Range("N218").Value = "`$/docena de matas"
$ws.Range("O218").Value = "Región de Coquimbo"
$ws.Range("P218").Value = 1251
$ws.Range("Q218").Value = 6
$ws.Range("R218").Value = "Hortaliza"

# Row 219 - new "Segunda" quality entry.
$ws.Range("A219").Value = 9
$ws.Range("B219").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value = "Metropolitana"
$ws.Range("D219").Value = 44610
$ws.Range("E219").Value = 13
$ws.Range("F219").Value = 100112017
$ws.Range("G219").Value = "Apio"
$ws.Range("H219").Value = "Americana (o)"
$ws.Range("I219").Value = "Segunda"
$ws.Range("J219").Value = 34
$ws.Range("K219").Value = 6000
$ws.Range("L219").Value = 6000
$ws.Range("M219").Value = 6000
$ws.Range("N219").Value = "`$/docena de matas"
$ws.Range("O219").Value = "Región de Coquimbo"
$ws.Range("P219").Value = 1000
$ws.Range("Q219").Value = 6
$ws.Range("R219").Value = "Hortaliza"
